# Sign off the timesheet: fill in supervisor name, supervisor signature
# initials and the sign-off date (as of 28/02/2014).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name (merged G6:I6)
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor Signature (merged A27:C27) + Date (merged D27:E27)
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = (Get-Date -Year 2014 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)

# Match the date cell's number format to the existing signed date cell (D25)
$ws.Range("D27").NumberFormat = $ws.Range("D25").NumberFormat

# Update the active selection / view to mirror the sign-off row.
$ws.Range("D27:E27").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
